# Add Stock_Min / Stock_Max columns and update some inventory quantities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (H1, I1) ---
$ws.Range("H1").Value = "Stock_Min"
$ws.Range("I1").Value = "Stock_Max"

# Match the existing header style (bold, bordered, centered) used by A1:G1
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

# --- Update existing Quantite values (column C) ---
$ws.Range("C2").Value = 30
$ws.Range("C4").Value = 9
$ws.Range("C5").Value = 1234
$ws.Range("C6").Value = 42
$ws.Range("C7").Value = 9

# --- Populate new Stock_Min (H) / Stock_Max (I) columns ---
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 50

$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 50

$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 30

$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 1200

$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 45

$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 100
